$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "Global MBA 교육 방식 – Essay 기반 평가란?"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/global-mba-essay-evaluation/#utm_source=rss&utm_medium=rss&utm_campaign=global-mba-essay-evaluation"

$ws.Range("D27").Value = "루다, 눈을 뜨다! 포토챗 베타의 멀티모달 기술 소개"
$ws.Range("E27").Value = "https://tech.scatterlab.co.kr/photochat-beta/"

$ws.Range("D50").Value = "디랙 벨트"
$ws.Range("E50").Value = "http://incredible.egloos.com/7573375"

$ws.Range("D51").Value = "데이터를 주고 받을 때 많이 사용되는 JSON 정리"
$ws.Range("E51").Value = "https://bskyvision.com/entry/%EB%8D%B0%EC%9D%B4%ED%84%B0%EB%A5%BC-%EC%A3%BC%EA%B3%A0-%EB%B0%9B%EC%9D%84-%EB%95%8C-%EB%A7%8E%EC%9D%B4-%EC%82%AC%EC%9A%A9%EB%90%98%EB%8A%94-JSON-%EC%A0%95%EB%A6%AC"
